$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0713
$ws.Range("E2").Value = -0.307
$ws.Range("G2").Value = 0.03741626794258374
$ws.Range("H2").Value = 0.03110047846889952
$ws.Range("I2").Value = 0.002631578947368421
$ws.Range("J2").Value = 0.001315789473684211
$ws.Range("K2").Value = 0.706
$ws.Range("L2").Value = 0.003377990430622009
$ws.Range("M2").Value = 1.3
$ws.Range("N2").Value = 0.02742616033755274
$ws.Range("O2").Value = 1.841359773371105
$ws.Range("P2").Value = 1.3
$ws.Range("Q2").Value = 0.02742616033755274
$ws.Range("R2").Value = 1.841359773371105
$ws.Range("U2").Value = 72.5
$ws.Range("V2").Value = 1.529535864978903
$ws.Range("W2").Value = 0.00815242494226328
$ws.Range("X2").Value = 0.07657972321424952
$ws.Range("Y2").Value = -0.06842729827198624
$ws.Range("Z2").Value = 13.66013071895425
$ws.Range("AA2").Value = 0.01797385620915033
$ws.Range("AB2").Value = 0.07652357922805854
$ws.Range("AC2").Value = -0.05854972301890821
$ws.Range("AD2").Value = 0.057
$ws.Range("AF2").Value = 0.057
$ws.Range("AG2").Value = -72.443
$ws.Range("AH2").Value = 0.00120108730008218
$ws.Range("AI2").Value = 0.0006212060115304554
$ws.Range("AJ2").Value = 2.892744479495268
$ws.Range("AK2").Value = -3.761904761904761
$ws.Range("AM2").Value = -1.05
$ws.Range("AN2").Value = 0.0103448275862069
$ws.Range("AP2").Value = -13.1475499092559
$ws.Range("AQ2").Value = -0.5238095238095238

$ws.Range("D3").Value = 0.0713
$ws.Range("E3").Value = -0.307
$ws.Range("G3").Value = 0.03741626794258374
$ws.Range("H3").Value = 0.03110047846889952
$ws.Range("I3").Value = 0.002631578947368421
$ws.Range("J3").Value = 0.001315789473684211
$ws.Range("K3").Value = 0.706
$ws.Range("L3").Value = 0.003377990430622009
$ws.Range("M3").Value = 1.3
$ws.Range("N3").Value = 0.02742616033755274
$ws.Range("O3").Value = 1.841359773371105
$ws.Range("P3").Value = 1.3
$ws.Range("Q3").Value = 0.02742616033755274
$ws.Range("R3").Value = 1.841359773371105
$ws.Range("U3").Value = 72.5
$ws.Range("V3").Value = 1.529535864978903
$ws.Range("W3").Value = 0.00815242494226328
$ws.Range("X3").Value = 0.07657972321424952
$ws.Range("Y3").Value = -0.06842729827198624
$ws.Range("Z3").Value = 13.66013071895425
$ws.Range("AA3").Value = 0.01797385620915033
$ws.Range("AB3").Value = 0.07652357922805854
$ws.Range("AC3").Value = -0.05854972301890821
$ws.Range("AD3").Value = 0.057
$ws.Range("AF3").Value = 0.057
$ws.Range("AG3").Value = -72.443
$ws.Range("AH3").Value = 0.00120108730008218
$ws.Range("AI3").Value = 0.0006212060115304554
$ws.Range("AJ3").Value = 2.892744479495268
$ws.Range("AK3").Value = -3.761904761904761
$ws.Range("AM3").Value = -1.05
$ws.Range("AN3").Value = 0.0103448275862069
$ws.Range("AP3").Value = -13.1475499092559
$ws.Range("AQ3").Value = -0.5238095238095238

